$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.678.18'
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").Value = '2.434.23'
$ws.Range("E3").Value = '  -1.41%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''505.91'
$ws.Range("E5").Value = '  -2.50%  '
$ws.Range("D6").Value = '''128.88'
$ws.Range("E6").Value = '  -2.61%  '
$ws.Range("D7").Value = '''0.998'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  -1.36%  '
$ws.Range("D9").Value = '2.447.37'
$ws.Range("E9").Value = '  -0.95%  '
$ws.Range("E10").Value = '  -0.18%  '
$ws.Range("D11").Value = '''0.0953'
$ws.Range("E11").Value = '  -3.95%  '
$ws.Range("D12").Value = '''5.17'
$ws.Range("E12").Value = '  -3.95%  '
$ws.Range("D13").Value = '''0.329'
$ws.Range("E13").Value = '  -3.66%  '
$ws.Range("D14").Value = '2.866.28'
$ws.Range("E14").Value = '  -1.40%  '
$ws.Range("D15").Value = '57.599.72'
$ws.Range("E15").Value = '  -0.75%  '
$ws.Range("D16").Value = '''21.86'
$ws.Range("E16").Value = '  -0.98%  '
$ws.Range("E17").Value = '  -2.71%  '
$ws.Range("D18").Value = '2.442.00'
$ws.Range("E18").Value = '  -1.01%  '
$ws.Range("D19").Value = '''10.45'
$ws.Range("E19").Value = '  -3.74%  '
$ws.Range("D20").Value = '''4.10'
$ws.Range("E20").Value = '  -1.84%  '
$ws.Range("D21").Value = '''314.68'
$ws.Range("E21").Value = '  -1.78%  '
$ws.Range("D22").Value = '''0.999'
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("E23").Value = '  -1.30%  '
$ws.Range("D24").Value = '''63.46'
$ws.Range("E24").Value = '  -1.11%  '
$ws.Range("D25").Value = '''0.406'
$ws.Range("E25").Value = '  -0.73%  '
$ws.Range("D26").Value = '''0.996'
$ws.Range("E26").Value = '  -0.27%  '
$ws.Range("E27").Value = '  -1.16%  '
$ws.Range("D28").Value = '''7.24'
$ws.Range("E28").Value = '  -2.04%  '
$ws.Range("D29").Value = '''169.37'
$ws.Range("E29").Value = '  +2.58%  '
$ws.Range("D30").Value = '0.0₃0724'
$ws.Range("E30").Value = '  -3.46%  '
$ws.Range("D31").Value = '''6.22'
$ws.Range("E31").Value = '  -3.16%  '
$ws.Range("D32").Value = '''1.65'
$ws.Range("E32").Value = '  -2.90%  '
$ws.Range("D33").Value = '''1.14'
$ws.Range("E33").Value = '  +0.67%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("E36").Value = '  -2.07%  '
$ws.Range("D37").Value = '''1.27'
$ws.Range("E37").Value = '  -5.20%  '
$ws.Range("D38").Value = '''3.92'
$ws.Range("E38").Value = '  -1.92%  '
$ws.Range("D39").Value = '''36.32'
$ws.Range("E39").Value = '  -0.48%  '
$ws.Range("E40").Value = '  -2.11%  '
$ws.Range("D41").Value = '''0.760'
$ws.Range("E41").Value = '  -4.57%  '
$ws.Range("D42").Value = '''271.12'
$ws.Range("E42").Value = '  -1.61%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '''3.37'
$ws.Range("E43").Value = '  -3.30%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '''4.97'
$ws.Range("E44").Value = '  -0.68%  '
$ws.Range("E45").Value = '  -2.16%  '
$ws.Range("D46").Value = '''0.0909'
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").Value = '''119.71'
$ws.Range("E47").Value = '  -5.54%  '
$ws.Range("D48").Value = '''0.0484'
$ws.Range("E48").Value = '  -1.61%  '
$ws.Range("D49").Value = '''17.16'
$ws.Range("E49").Value = '  -3.55%  '
$ws.Range("E50").Value = '  -2.29%  '
$ws.Range("D51").Value = '''16.61'
$ws.Range("E51").Value = '  -3.13%  '
